$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap full row content (columns B:AD) between row 19 and row 20, keep column A (index) fixed
$ws.Range("B19").Value = 7032917
$ws.Range("B20").Value = 7032914
$ws.Range("E19").Value = "FK Backa Topola"
$ws.Range("E20").Value = "FK Vozdovac"
$ws.Range("F19").Value = "FK Radnicki 1923"
$ws.Range("F20").Value = "FK Radnik Surdulica"
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 1
$ws.Range("J19").Value = 0
$ws.Range("J20").Value = 1
$ws.Range("K19").Value = "H"
$ws.Range("K20").Value = "D"
$ws.Range("L19").Value = 1.5
$ws.Range("L20").Value = 2.2
$ws.Range("M19").Value = 3.75
$ws.Range("M20").Value = 3.1
$ws.Range("N19").Value = 6.5
$ws.Range("N20").Value = 3.2
$ws.Range("O19").Value = 1.444
$ws.Range("O20").Value = 2.05
$ws.Range("P19").Value = 4
$ws.Range("P20").Value = 3.1
$ws.Range("Q19").Value = 6.5
$ws.Range("Q20").Value = 3.5
$ws.Range("R19").Value = -1.25
$ws.Range("R20").Value = -0.25
$ws.Range("S19").Value = 1.975
$ws.Range("S20").Value = 1.75
$ws.Range("T19").Value = 1.825
$ws.Range("T20").Value = 2.05
$ws.Range("U19").Value = 2.75
$ws.Range("U20").Value = 2
$ws.Range("V19").Value = 1.95
$ws.Range("V20").Value = 1.775
$ws.Range("W19").Value = 1.85
$ws.Range("W20").Value = 2.025
$ws.Range("X19").Value = 0.444
$ws.Range("X20").Value = -1
$ws.Range("Y19").Value = -1
$ws.Range("Y20").Value = 2.1
$ws.Range("AB19").Value = 0.4125
$ws.Range("AB20").Value = 0.5249999999999999
$ws.Range("AC19").Value = -1
$ws.Range("AC20").Value = 0
$ws.Range("AD19").Value = 0.8500000000000001
$ws.Range("AD20").Value = 0

# Swap full row content (columns B:AD) between row 31 and row 32, keep column A (index) fixed
$ws.Range("B31").Value = 6979431
$ws.Range("B32").Value = 6979427
$ws.Range("E31").Value = "FK Napredak"
$ws.Range("E32").Value = "Javor Ivanjica"
$ws.Range("F31").Value = "Mladost Lucani"
$ws.Range("F32").Value = "FK Backa Topola"
$ws.Range("H31").Value = 1
$ws.Range("H32").Value = 3
$ws.Range("J31").Value = 0
$ws.Range("J32").Value = 3
$ws.Range("L31").Value = 1.7
$ws.Range("L32").Value = 5.25
$ws.Range("M31").Value = 3.75
$ws.Range("M32").Value = 4
$ws.Range("N31").Value = 4.5
$ws.Range("N32").Value = 1.533
$ws.Range("O31").Value = 1.909
$ws.Range("O32").Value = 4
$ws.Range("P31").Value = 3.1
$ws.Range("P32").Value = 3.3
$ws.Range("Q31").Value = 4
$ws.Range("Q32").Value = 1.85
$ws.Range("R31").Value = -0.5
$ws.Range("R32").Value = 0.5
$ws.Range("S31").Value = 1.975
$ws.Range("S32").Value = 1.875
$ws.Range("T31").Value = 1.825
$ws.Range("T32").Value = 1.925
$ws.Range("U31").Value = 2
$ws.Range("U32").Value = 2.25
$ws.Range("V31").Value = 1.875
$ws.Range("V32").Value = 1.825
$ws.Range("W31").Value = 1.925
$ws.Range("W32").Value = 1.975
$ws.Range("Z31").Value = 3
$ws.Range("Z32").Value = 0.8500000000000001
$ws.Range("AB31").Value = 0.825
$ws.Range("AB32").Value = 0.925
$ws.Range("AC31").Value = -1
$ws.Range("AC32").Value = 0.825
$ws.Range("AD31").Value = 0.925
$ws.Range("AD32").Value = -1

# Swap full row content (columns B:AD) between row 77 and row 78, keep column A (index) fixed
$ws.Range("B77").Value = 6979476
$ws.Range("B78").Value = 6979475
$ws.Range("E77").Value = "FK Backa Topola"
$ws.Range("E78").Value = "FK Cukaricki"
$ws.Range("F77").Value = "Spartak Subotica"
$ws.Range("F78").Value = "FK Vozdovac"
$ws.Range("G77").Value = 0
$ws.Range("G78").Value = 1
$ws.Range("H77").Value = 2
$ws.Range("H78").Value = 1
$ws.Range("J77").Value = 0
$ws.Range("J78").Value = 1
$ws.Range("K77").Value = "A"
$ws.Range("K78").Value = "D"
$ws.Range("L77").Value = 1.4
$ws.Range("L78").Value = 1.533
$ws.Range("M77").Value = 4.2
$ws.Range("M78").Value = 3.6
$ws.Range("N77").Value = 6
$ws.Range("N78").Value = 5.5
$ws.Range("O77").Value = 1.363
$ws.Range("O78").Value = 1.444
$ws.Range("P77").Value = 4.5
$ws.Range("P78").Value = 3.75
$ws.Range("Q77").Value = 6.5
$ws.Range("Q78").Value = 7
$ws.Range("S77").Value = 1.85
$ws.Range("S78").Value = 2
$ws.Range("T77").Value = 1.95
$ws.Range("T78").Value = 1.8
$ws.Range("U77").Value = 2.75
$ws.Range("U78").Value = 2.5
$ws.Range("V77").Value = 1.85
$ws.Range("V78").Value = 1.95
$ws.Range("W77").Value = 1.95
$ws.Range("W78").Value = 1.85
$ws.Range("Y77").Value = -1
$ws.Range("Y78").Value = 2.75
$ws.Range("Z77").Value = 5.5
$ws.Range("Z78").Value = -1
$ws.Range("AB77").Value = 0.95
$ws.Range("AB78").Value = 0.8
$ws.Range("AD77").Value = 0.95
$ws.Range("AD78").Value = 0.8500000000000001

# Swap full row content (columns B:AD) between row 90 and row 91, keep column A (index) fixed
$ws.Range("B90").Value = 6978747
$ws.Range("B91").Value = 6979491
$ws.Range("E90").Value = "IMT Novi Belgrade"
$ws.Range("E91").Value = "Radnicki Nis"
$ws.Range("F90").Value = "Red Star Belgrade"
$ws.Range("F91").Value = "Spartak Subotica"
$ws.Range("H90").Value = 2
$ws.Range("H91").Value = 1
$ws.Range("J90").Value = 2
$ws.Range("J91").Value = 1
$ws.Range("K90").Value = "A"
$ws.Range("K91").Value = "D"
$ws.Range("L90").Value = 8
$ws.Range("L91").Value = 1.95
$ws.Range("M90").Value = 5.25
$ws.Range("M91").Value = 3.25
$ws.Range("N90").Value = 1.285
$ws.Range("N91").Value = 3.7
$ws.Range("O90").Value = 15
$ws.Range("O91").Value = 1.65
$ws.Range("P90").Value = 7.5
$ws.Range("P91").Value = 3.5
$ws.Range("Q90").Value = 1.125
$ws.Range("Q91").Value = 5
$ws.Range("R90").Value = 2.25
$ws.Range("R91").Value = -0.75
$ws.Range("S90").Value = 1.975
$ws.Range("S91").Value = 1.825
$ws.Range("T90").Value = 1.825
$ws.Range("T91").Value = 1.975
$ws.Range("U90").Value = 3.5
$ws.Range("U91").Value = 2.5
$ws.Range("V90").Value = 1.825
$ws.Range("V91").Value = 2
$ws.Range("W90").Value = 1.975
$ws.Range("W91").Value = 1.8
$ws.Range("Y90").Value = -1
$ws.Range("Y91").Value = 2.5
$ws.Range("Z90").Value = 0.125
$ws.Range("Z91").Value = -1
$ws.Range("AA90").Value = 0.9750000000000001
$ws.Range("AA91").Value = -1
$ws.Range("AB90").Value = -1
$ws.Range("AB91").Value = 0.9750000000000001
$ws.Range("AD90").Value = 0.9750000000000001
$ws.Range("AD91").Value = 0.8

# Swap full row content (columns B:AD) between row 167 and row 168, keep column A (index) fixed
$ws.Range("B167").Value = 6979545
$ws.Range("B168").Value = 6979547
$ws.Range("E167").Value = "Radnicki Nis"
$ws.Range("E168").Value = "FK Backa Topola"
$ws.Range("F167").Value = "Javor Ivanjica"
$ws.Range("F168").Value = "FK Radnik Surdulica"
$ws.Range("I167").Value = 1
$ws.Range("I168").Value = 0
$ws.Range("L167").Value = 2
$ws.Range("L168").Value = 1.333
$ws.Range("M167").Value = 3.25
$ws.Range("M168").Value = 4.333
$ws.Range("N167").Value = 3.25
$ws.Range("N168").Value = 7.5
$ws.Range("O167").Value = 1.727
$ws.Range("O168").Value = 1.25
$ws.Range("P167").Value = 3.3
$ws.Range("P168").Value = 4.75
$ws.Range("Q167").Value = 4.2
$ws.Range("Q168").Value = 10
$ws.Range("R167").Value = -0.5
$ws.Range("R168").Value = -1.5
$ws.Range("S167").Value = 1.825
$ws.Range("S168").Value = 1.85
$ws.Range("T167").Value = 1.975
$ws.Range("T168").Value = 1.95
$ws.Range("U167").Value = 2.25
$ws.Range("U168").Value = 2.5
$ws.Range("V167").Value = 1.825
$ws.Range("V168").Value = 1.8
$ws.Range("W167").Value = 1.975
$ws.Range("W168").Value = 2
$ws.Range("X167").Value = 0.7270000000000001
$ws.Range("X168").Value = 0.25
$ws.Range("AA167").Value = 0.825
$ws.Range("AA168").Value = -1
$ws.Range("AB167").Value = -1
$ws.Range("AB168").Value = 0.95
$ws.Range("AD167").Value = 0.9750000000000001
$ws.Range("AD168").Value = 1
